$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header in K1 matching the bold/centered style used by the other headers (row 1)
$ws.Range("K1").Value = "PhylogenySorting"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108  # xlCenter

# Add the new data value in K4
$ws.Range("K4").Value = "T45"

# Match the column widths recorded in the workbook after the edit
# (values chosen so the engine's internal rounding lands on the saved width)
$ws.Columns.Item(10).ColumnWidth = 11.833333333333332
$ws.Columns.Item(11).ColumnWidth = 16.833333333333336

# Leave the selection where it ended up when the file was last saved
$ws.Range("G8").Select()
